$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '62.301.88'
$ws.Range('E2').Value = '  -3.22%  '
Set-TextValue 'D3' '3.377.34'
$ws.Range('E3').Value = '  -3.80%  '
Set-TextValue 'D4' '0.999'
$ws.Range('E4').Value = '  -0.04%  '
Set-TextValue 'D5' '568.05'
$ws.Range('E5').Value = '  -3.86%  '
Set-TextValue 'D6' '124.79'
$ws.Range('E6').Value = '  -7.27%  '
$ws.Range('E7').Value = '  -0.01%  '
Set-TextValue 'D8' '3.379.02'
$ws.Range('E8').Value = '  -3.71%  '
Set-TextValue 'D9' '0.471'
$ws.Range('E9').Value = '  -3.40%  '
$ws.Range('E10').Value = '  -5.26%  '
$ws.Range('E11').Value = '  -4.63%  '
Set-TextValue 'D12' '0.373'
$ws.Range('E12').Value = '  -4.16%  '
Set-TextValue 'D13' '3.942.10'
$ws.Range('E13').Value = '  -4.02%  '
$ws.Range('E14').Value = '  -1.09%  '
Set-TextValue 'D15' '3.367.76'
$ws.Range('E15').Value = '  -3.99%  '
Set-TextValue 'D16' '0.0000170'
$ws.Range('E16').Value = '  -6.04%  '
Set-TextValue 'D17' '62.298.35'
$ws.Range('E17').Value = '  -3.19%  '
Set-TextValue 'D18' '24.36'
$ws.Range('E18').Value = '  -5.31%  '
Set-TextValue 'D19' '9.17'
$ws.Range('E19').Value = '  -8.32%  '
Set-TextValue 'D20' '5.61'
$ws.Range('E20').Value = '  -2.64%  '
Set-TextValue 'D21' '13.05'
$ws.Range('E21').Value = '  -3.93%  '
Set-TextValue 'D22' '369.99'
$ws.Range('E22').Value = '  -6.22%  '
Set-TextValue 'D23' '0.552'
$ws.Range('E23').Value = '  -4.43%  '
Set-TextValue 'D24' '3.504.67'
$ws.Range('E24').Value = '  -3.97%  '
$ws.Range('E25').Value = '  -0.05%  '
Set-TextValue 'D26' '70.76'
Set-TextValue 'D27' '0.0000105'
$ws.Range('E27').Value = '  -10.95%  '
Set-TextValue 'D28' '0.996'
$ws.Range('E28').Value = '  -0.45%  '
Set-TextValue 'D29' '6.84'
$ws.Range('E29').Value = '  -7.31%  '
$ws.Range('E30').Value = '  -6.79%  '
Set-TextValue 'D31' '7.75'
$ws.Range('E31').Value = '  -6.10%  '
Set-TextValue 'D33' '3.399.49'
$ws.Range('E33').Value = '  -3.96%  '
$ws.Range('E34').Value = '  -6.00%  '
$ws.Range('E35').Value = '  -6.60%  '
Set-TextValue 'D36' '22.58'
$ws.Range('E36').Value = '  -3.40%  '
Set-TextValue 'D37' '5.08'
$ws.Range('E37').Value = '  -5.24%  '
Set-TextValue 'D38' '164.64'
$ws.Range('E38').Value = '  -1.64%  '
$ws.Range('E39').Value = '  -5.57%  '
Set-TextValue 'D40' '1.47'
$ws.Range('E40').Value = '  -5.67%  '
Set-TextValue 'D41' '0.0748'
$ws.Range('E41').Value = '  -4.95%  '
Set-TextValue 'D42' '0.999'
$ws.Range('E42').Value = '  -0.15%  '
Set-TextValue 'D43' '0.763'
$ws.Range('E43').Value = '  -5.91%  '
Set-TextValue 'D44' '41.04'
$ws.Range('E44').Value = '  -2.53%  '
Set-TextValue 'D45' '4.22'
$ws.Range('E45').Value = '  -4.94%  '
$ws.Range('E46').Value = '  -7.79%  '
Set-TextValue 'D47' '22.28'
$ws.Range('E47').Value = '  -11.13%  '
Set-TextValue 'D48' '1.06'
$ws.Range('E48').Value = '  -9.83%  '
Set-TextValue 'D49' '6.58'
$ws.Range('E49').Value = '  -3.32%  '
Set-TextValue 'D50' '2.227.00'
$ws.Range('E50').Value = '  -6.47%  '
Set-TextValue 'D51' '0.837'
$ws.Range('E51').Value = '  -6.54%  '
